# Daily attendance processing - normalize "Recorded By" (column G) ordering
# so that the "System" token is listed first among the comma-separated
# recorder names, for every row that currently has it listed elsewhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notmatch ",") { continue }

    $parts = $val -split ", "

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p -ieq "system") { $hasSystem = $true }
    }
    if (-not $hasSystem) { continue }

    $reversedParts = $parts[($parts.Count - 1)..0]
    $newVal = [string]::Join(", ", $reversedParts)

    if (-not $newVal.Equals($val)) {
        $cell.Value = $newVal
    }
}
